$d = $word.ActiveDocument
$r = $d.Range(0, 5)
Write-Host "r text: [$($r.Text)]"
$r.Paragraphs.LineSpacingRule = 0
